{"js": "// Apply the literature-review text edits described by the commit\n// \"more lit review again\": update the version-date stamp, fix a couple\n// of typos, and insert several additional citation/detail clauses into\n// existing bullet points and paragraphs. Every target string below is\n// the *entire* contents of a single run's <w:t> in the source document,\n// so a plain exact-text search-and-replace is sufficient and safe.\nconst replacements = [\n  [\"Version Date: Thu Aug 23 18:03:41 2018 -0700\", \"Version Date: Fri Aug 24 12:01:20 2018 -0700\"],\n  [\"Progress in the field of applied bioassessment is impeded by a lack of methods that are accessible and reproducible for the management community. Open science principles that seek to democratize science can address these challenges, yet widespread adoption in research has yet to gain traction for the devepment and appplication bioassessment methods. At the core of this philosophy is the concept that research should be reproducible and transparent, in addition to having long-term provenance through effective modes of data preservation and sharing. This review will introduce core open science concepts that have been advocated more generally in the ecological sciences and will emphasize how adoption can benefit bioassessment for both prescriptive condition assessments and proactive applications that inform planning activities. Examples from the state of California will be used to demonstrate effective adoption of open science principles through data stewardship, reproducible research, and engagement of stakeholders with multimedia applications. Technical, sociocultural, and institutional challenges for adopting open science will also be discussed, including practical approaches for overcoming these hurdles in bioassessment applications.\", \"Progress in the field of applied bioassessment is impeded by a lack of methods that are accessible and reproducible for the management community. Open science principles that seek to democratize science can address these challenges, yet widespread adoption in research has yet to gain traction for the development and appplication bioassessment methods. At the core of this philosophy is the concept that research should be reproducible and transparent, in addition to having long-term provenance through effective modes of data preservation and sharing. This review will introduce core open science concepts that have been advocated more generally in the ecological sciences and will emphasize how adoption can benefit bioassessment for both prescriptive condition assessments and proactive applications that inform planning activities. Examples from the state of California will be used to demonstrate effective adoption of open science principles through data stewardship, reproducible research, and engagement of stakeholders with multimedia applications. Technical, sociocultural, and institutional challenges for adopting open science will also be discussed, including practical approaches for overcoming these hurdles in bioassessment applications.\"],\n  [\"Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in overcoming technical challenges for developing accurate and interpretable bioassessment methods. However, the widespread use of bioassessment data by managers and stakeholders is severely imbalanced relative to the number of indices that have been developed . Existing methods often lack transparency , require specialized training to implement, and are not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.\", \"Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in overcoming technical challenges for developing accurate and interpretable assessment methods that rely on biological organisms as sentinels of environmental condition. However, the widespread use of bioassessment data by managers and stakeholders is severely imbalanced relative to the number of indices that have been developed. Existing methods often lack transparency, require specialized training to implement, and are not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.\"],\n  [\"Bioassessment applications are usually mandated by legislation and informs management of resources through condition assessment\", \"Bioassessment applications are usually mandated by legislation and informs management of resources through condition assessment - CWA in US, WFD in Europe\"],\n  [\"Proliferation of methods with lack of transparency - review of index coverage in US, internationally\", \"Proliferation of methods - review of index coverage in US, internationally, Birk et al. 2012 describe nearly 300 methods that have been developed in Europe, Nichols et al. 2016 describe national down-scaling of bioassessment applications in Australia as a negative in absence of coordinated federal assessment networks, Kelly et al. 2016 describe issues of redundancy in ecological assessment of lakes regarding use of multiple taxa, are they all necessary?\"],\n  [\"Lack of access to index calibration/validation data, information that is often collected through public funds although often treated as proprietary\", \"Lack of access to index calibration/validation data, information that is often collected through public funds although often treated as proprietary, data are not always discoverable (Hering et al. 2010 describe this issue in a ten-year assessment of WFD)\"],\n  [\"of the ecological sciences, i.e., the carefully collected observational data meant to address specific research questions. Scientists in the long-tail are potentially more relucant to adopt open-science because of the perception of less benefit to making the data open. This suggests that bioassessment datasets and associated methods are inherently more likely to benefit from openness because more widespread appeal. Conversely, the long-tail datasets individually may not have broad relevance but collectively could serve larger purposes (Jenny\\u2019s project?)\", \"of the ecological sciences, i.e., the carefully collected observational data meant to address specific research questions. Scientists in the long-tail are potentially more relucant to adopt open-science because of the perception of less benefit to making the data open. This suggests that bioassessment datasets and associated methods are inherently more likely to benefit from openness because more widespread appeal. Conversely, the long-tail datasets individually may not have broad relevance but collectively could serve larger purposes, some countries have abandonated national-scale coordinated monitoring efforts in favor localized sampling (Nichols et al. 2016)\"],\n  [\"Call to implement now - field is transitioning to molecular approaches where information acquisition will be orders of magnitude greater than traditional taxonomic-based approaches. Data acquisition and management will require systematic methods for documenting, cataloging, and sharing information - start now. Use of online eDNA archives have been established.\", \"Call to implement now - field is transitioning to molecular approaches where information acquisition will be orders of magnitude greater than traditional taxonomic-based approaches. Data acquisition and management will require systematic methods for documenting, cataloging, and sharing information - start now. Use of online eDNA archives have been established. Baird and Hajibabaei 2012 describe the bioassessment paradigm with molecular approaches\"]\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText.substring(0, 60));\n  }\n\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the literature-review text edits described by the commit\n# \"more lit review again\": update the version-date stamp, fix a couple\n# of typos, and insert several additional citation/detail clauses into\n# existing bullet points and paragraphs. Every target string below is\n# the *entire* contents of a single run, so an exact-text Find/Replace\n# (no wildcards) on $d.Content is sufficient and safe.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($OldText, $NewText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $OldText\n    $find.Replacement.Text = $NewText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdFindContinue wrap (1), wdReplaceOne (1): replace only the single\n    # exact occurrence located by this search.\n    $found = $find.Execute($OldText, $false, $true, $false, $false, $false, $true, 1, $false, $NewText, 1)\n    if (-not $found) {\n        throw \"Replace-ExactText: text not found: $OldText\"\n    }\n}\n\nReplace-ExactText \"Version Date: Thu Aug 23 18:03:41 2018 -0700\" \"Version Date: Fri Aug 24 12:01:20 2018 -0700\"\nReplace-ExactText \"Progress in the field of applied bioassessment is impeded by a lack of methods that are accessible and reproducible for the management community. Open science principles that seek to democratize science can address these challenges, yet widespread adoption in research has yet to gain traction for the devepment and appplication bioassessment methods. At the core of this philosophy is the concept that research should be reproducible and transparent, in addition to having long-term provenance through effective modes of data preservation and sharing. This review will introduce core open science concepts that have been advocated more generally in the ecological sciences and will emphasize how adoption can benefit bioassessment for both prescriptive condition assessments and proactive applications that inform planning activities. Examples from the state of California will be used to demonstrate effective adoption of open science principles through data stewardship, reproducible research, and engagement of stakeholders with multimedia applications. Technical, sociocultural, and institutional challenges for adopting open science will also be discussed, including practical approaches for overcoming these hurdles in bioassessment applications.\" \"Progress in the field of applied bioassessment is impeded by a lack of methods that are accessible and reproducible for the management community. Open science principles that seek to democratize science can address these challenges, yet widespread adoption in research has yet to gain traction for the development and appplication bioassessment methods. At the core of this philosophy is the concept that research should be reproducible and transparent, in addition to having long-term provenance through effective modes of data preservation and sharing. This review will introduce core open science concepts that have been advocated more generally in the ecological sciences and will emphasize how adoption can benefit bioassessment for both prescriptive condition assessments and proactive applications that inform planning activities. Examples from the state of California will be used to demonstrate effective adoption of open science principles through data stewardship, reproducible research, and engagement of stakeholders with multimedia applications. Technical, sociocultural, and institutional challenges for adopting open science will also be discussed, including practical approaches for overcoming these hurdles in bioassessment applications.\"\nReplace-ExactText \"Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in overcoming technical challenges for developing accurate and interpretable bioassessment methods. However, the widespread use of bioassessment data by managers and stakeholders is severely imbalanced relative to the number of indices that have been developed . Existing methods often lack transparency , require specialized training to implement, and are not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.\" \"Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in overcoming technical challenges for developing accurate and interpretable assessment methods that rely on biological organisms as sentinels of environmental condition. However, the widespread use of bioassessment data by managers and stakeholders is severely imbalanced relative to the number of indices that have been developed. Existing methods often lack transparency, require specialized training to implement, and are not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.\"\nReplace-ExactText \"Bioassessment applications are usually mandated by legislation and informs management of resources through condition assessment\" \"Bioassessment applications are usually mandated by legislation and informs management of resources through condition assessment - CWA in US, WFD in Europe\"\nReplace-ExactText \"Proliferation of methods with lack of transparency - review of index coverage in US, internationally\" \"Proliferation of methods - review of index coverage in US, internationally, Birk et al. 2012 describe nearly 300 methods that have been developed in Europe, Nichols et al. 2016 describe national down-scaling of bioassessment applications in Australia as a negative in absence of coordinated federal assessment networks, Kelly et al. 2016 describe issues of redundancy in ecological assessment of lakes regarding use of multiple taxa, are they all necessary?\"\nReplace-ExactText \"Lack of access to index calibration/validation data, information that is often collected through public funds although often treated as proprietary\" \"Lack of access to index calibration/validation data, information that is often collected through public funds although often treated as proprietary, data are not always discoverable (Hering et al. 2010 describe this issue in a ten-year assessment of WFD)\"\nReplace-ExactText \"of the ecological sciences, i.e., the carefully collected observational data meant to address specific research questions. Scientists in the long-tail are potentially more relucant to adopt open-science because of the perception of less benefit to making the data open. This suggests that bioassessment datasets and associated methods are inherently more likely to benefit from openness because more widespread appeal. Conversely, the long-tail datasets individually may not have broad relevance but collectively could serve larger purposes (Jenny\u2019s project?)\" \"of the ecological sciences, i.e., the carefully collected observational data meant to address specific research questions. Scientists in the long-tail are potentially more relucant to adopt open-science because of the perception of less benefit to making the data open. This suggests that bioassessment datasets and associated methods are inherently more likely to benefit from openness because more widespread appeal. Conversely, the long-tail datasets individually may not have broad relevance but collectively could serve larger purposes, some countries have abandonated national-scale coordinated monitoring efforts in favor localized sampling (Nichols et al. 2016)\"\nReplace-ExactText \"Call to implement now - field is transitioning to molecular approaches where information acquisition will be orders of magnitude greater than traditional taxonomic-based approaches. Data acquisition and management will require systematic methods for documenting, cataloging, and sharing information - start now. Use of online eDNA archives have been established.\" \"Call to implement now - field is transitioning to molecular approaches where information acquisition will be orders of magnitude greater than traditional taxonomic-based approaches. Data acquisition and management will require systematic methods for documenting, cataloging, and sharing information - start now. Use of online eDNA archives have been established. Baird and Hajibabaei 2012 describe the bioassessment paradigm with molecular approaches\"\n"}
